# Stack.pptx edit script
# Commit: "changes related to the double ended queue and circular queue : VP"
#
# Summary of content changes (slide numbers per the open presentation):
#   Slide 1 (Stack overview):
#     - "ele Peek() ... inserted into the stack"  -> "ele Peek() ... last element "
#     - "ele[] Display() ..." run split into "ele" + "[] Display() ..."
#     - TextBox 50 (bullet list) shrinks height (autofit after text shortened)
#     - TextBox 13 ("top" pointer label) moves down/right to point at a
#       different array slot
#   Slide 2 (Push):
#     - "Check stack is overflow or not top == max -1" -> "... (top == max -1)"
#     - "top = 3"  -> "top = 2"
#   Slide 4 (Pop) and Slide 6 (Peek):
#     - "Check stack is underflow or not top == -1" -> "... (top == -1)"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1 : "Stack" overview slide
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# TextBox 50 holds the bulleted description / operations list.
$bulletBox = $s1.Shapes.Item(4)
$tr = $bulletBox.TextFrame.TextRange
$full = $tr.Text

# -- "ele Peek() ..." bullet: drop the trailing "inserted into the stack"
$oldPeek = "ele Peek() " + [char]0x2013 + " returns last element inserted into the stack"
$newPeek = "ele Peek() " + [char]0x2013 + " returns last element "
$idx = $full.IndexOf($oldPeek)
if ($idx -ge 0) {
    $rng = $tr.Characters($idx + 1, $oldPeek.Length)
    $rng.Text = $newPeek
}

# -- "ele[] Display() ..." bullet: split "ele" into its own run (so it can be
#    flagged independently, matching the authored edit) while leaving the
#    wording itself unchanged.
$full = $bulletBox.TextFrame.TextRange.Text
$eleDisplay = "ele[] Display()"
$idx2 = $full.IndexOf($eleDisplay)
if ($idx2 -ge 0) {
    $eleRng = $tr.Characters($idx2 + 1, 3)
    $eleRng.Text = "ele"
}

# Resize the bullet textbox (it auto-shrinks once the Peek bullet above got
# shorter).
$bulletBox.Height = 298.0828400456693

# "top" pointer label textbox moves to line up with a different array slot.
$topLabel = $s1.Shapes.Item(7)
$topLabel.Left = 644.4642639685039
$topLabel.Top = 243.9344101488189

# ---------------------------------------------------------------------------
# Slide 2 : "Push"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$pushBox = $s2.Shapes.Item(4)
$pushText = $pushBox.TextFrame.TextRange.Text
$oldOverflow = "Check stack is overflow or not top == max -1"
$newOverflow = "Check stack is overflow (top == max -1)"
$oidx = $pushText.IndexOf($oldOverflow)
if ($oidx -ge 0) {
    $orng = $pushBox.TextFrame.TextRange.Characters($oidx + 1, $oldOverflow.Length)
    $orng.Text = $newOverflow
}

$topThree = $s2.Shapes.Item(13)
$topThree.TextFrame.TextRange.Text = "top = 2"

# ---------------------------------------------------------------------------
# Slide 4 : "Pop"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

$popBox = $s4.Shapes.Item(4)
$popText = $popBox.TextFrame.TextRange.Text
$oldUnderflow = "Check stack is underflow or not top == -1"
$newUnderflow = "Check stack is underflow (top == -1)"
$uidx = $popText.IndexOf($oldUnderflow)
if ($uidx -ge 0) {
    $urng = $popBox.TextFrame.TextRange.Characters($uidx + 1, $oldUnderflow.Length)
    $urng.Text = $newUnderflow
}

# ---------------------------------------------------------------------------
# Slide 6 : "Peek"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

$peekBox = $s6.Shapes.Item(4)
$peekText = $peekBox.TextFrame.TextRange.Text
$uidx2 = $peekText.IndexOf($oldUnderflow)
if ($uidx2 -ge 0) {
    $urng2 = $peekBox.TextFrame.TextRange.Characters($uidx2 + 1, $oldUnderflow.Length)
    $urng2.Text = $newUnderflow
}
